$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a "numEpisodes" column right after "maxNumTries" (column D), pushing
# the Qlearn/SARSA mean+std headers one column to the right and dropping the
# trailing "SARSA0.75 - STD" column (the scale-up-factor column takes its
# former slot in the 12-column layout).
$ws.Range("L1").Value = "SARSA0.75 - Mean"
$ws.Range("K1").Value = "SARSA0.50 - STD"
$ws.Range("J1").Value = "SARSA0.5 - Mean"
$ws.Range("I1").Value = "SARSA0.25 - STD"
$ws.Range("H1").Value = "SARSA0.25 - Mean"
$ws.Range("G1").Value = "Qlearn - STD"
$ws.Range("F1").Value = "Qlearn - Mean"
$ws.Range("E1").Value = "numEpisodes"

# Update the active selection to F5 (as recorded in the post-edit sheet view)
$ws.Range("F5").Select()
